# Burndown chart update: record daily progress for tasks 1 and 2, and
# correct the initial estimate for task 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Task 2 (row 7) initial estimate corrected from 7 to 10
$ws.Range("D7").Value = 10

# Task 1 (row 6) progress logged on Day 8 / Day 9
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 3

# Task 2 (row 7) progress logged on Day 5 / Day 6 / Day 7
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 4

# Reflect the current view/selection as left by the author
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("K8").Select()
